$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell B1 - same value/style treatment as A1
$ws.Range("B1").Value = "Melhores DPS em Míticas"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

$data = @(
    "Shadow Priest (S Tier)",
    "Fire Mage (A Tier)",
    "Augmentation Evoker(A Tier)",
    "Frost Mage (A Tier)",
    "Outlaw Rogue (A Tier)",
    "Destruction Warlock (A Tier)",
    "Devastation Evoker (B Tier)",
    "Balance Druid (B Tier)",
    "Demonology Warlock (B Tier)",
    "Arms Warrior (B Tier)",
    "Windwalker Monk (B Tier)",
    "Beast Mastery Hunter (B Tier)",
    "Demonology Warlock (B Tier)",
    "Subtlety Rogue (B Tier)",
    "Retribution Paladin (B Tier)",
    "Assassination Rogue (C Tier)",
    "Marksmanship Hunter (C Tier)",
    "Fury Warrior (C Tier)",
    "Unholy Death Knight (C Tier)",
    "Elemental Shaman (C Tier)",
    "Frost Death Knight (C Tier)",
    "Arcane Mage (C Tier)",
    "Assassination Rogue (C Tier)",
    "Enhancement Shaman (C Tier)",
    "Havoc Demon Hunter (C Tier)",
    "Affliction Warlock (D Tier)",
    "Feral Druid (D Tier)",
    "Survival Hunter (D Tier)"
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i]
}
